# Adds several new glossary definitions to the "Nomenclature" section, and
# relocates the "_GoBack" bookmark from the end of the "Recursion" paragraph
# to the end of the newly-added "Integrated Development Environment (IDE)"
# paragraph (matching the target diff).

$d = $word.ActiveDocument

# Remove the original "_GoBack" bookmark (it currently sits at the end of
# the "Recursion" paragraph) up front so that later, when a fresh "_GoBack"
# bookmark is created for the new "Integrated Development Environment
# (IDE)" paragraph, there is no name collision.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

function Find-Paragraph {
    param($doc, $searchText)
    $rng = $doc.Content
    $null = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return $rng.Paragraphs(1)
}

function Get-ParaIndex {
    param($doc, $para)
    $idx = 0
    foreach ($pp in $doc.Paragraphs) {
        $idx = $idx + 1
        if ($pp.Range.Start -eq $para.Range.Start) {
            return $idx
        }
    }
    return -1
}

function Insert-TermBefore {
    # Inserts a new "Term: definition" paragraph immediately before $anchorPara.
    # $defParts is an array of definition-text chunks (all rendered non-bold);
    # pass an empty array to create a term-only paragraph.
    param($doc, $anchorPara, $term, $defParts)

    $idx = Get-ParaIndex $doc $anchorPara
    $r = $anchorPara.Range
    $r.Collapse(1)
    $r.InsertParagraphBefore()

    $newPara = $doc.Paragraphs($idx)
    $nr = $newPara.Range
    $nr.Text = $term
    $boldRange = $doc.Range($nr.Start, $nr.Start + $term.Length)
    $boldRange.Font.Bold = 1

    foreach ($def in $defParts) {
        $insPoint = $newPara.Range.End - 1
        $defRange = $doc.Range($insPoint, $insPoint)
        $defRange.InsertAfter($def)
        $defRange.Font.Bold = 0
    }

    return $newPara
}

function Insert-TermAfter {
    # Inserts a new "Term: definition" paragraph immediately after $anchorPara.
    param($doc, $anchorPara, $term, $defParts)

    $idx = Get-ParaIndex $doc $anchorPara
    $r = $anchorPara.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()

    $newPara = $doc.Paragraphs($idx + 1)
    $nr = $newPara.Range
    $nr.Text = $term
    $boldRange = $doc.Range($nr.Start, $nr.Start + $term.Length)
    $boldRange.Font.Bold = 1

    foreach ($def in $defParts) {
        $insPoint = $newPara.Range.End - 1
        $defRange = $doc.Range($insPoint, $insPoint)
        $defRange.InsertAfter($def)
        $defRange.Font.Bold = 0
    }

    return $newPara
}

function Add-BookmarkAtParagraphEnd {
    # Adds a zero-length bookmark named $bmName at the very end of $para's
    # text (i.e. right before the paragraph mark), by replacing the
    # paragraph's content (minus the trailing mark) via InsertXML, since
    # Bookmarks.Add placed exactly on a paragraph-mark position is unreliable.
    param($doc, $para, $bmName)

    $fullRange = $doc.Range($para.Range.Start, $para.Range.End - 1)
    $innerXml = $fullRange.WordOpenXML
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Integrated Development Environment (IDE)</w:t></w:r>' +
        '<w:r><w:t>: a source code editor with automation tools</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="' + $bmName + '"/><w:bookmarkEnd w:id="0"/></w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $fullRange.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 1) "Application Program Interface (API)" before "Constant"
# ---------------------------------------------------------------------------
$constantPara = Find-Paragraph $d "Constant:"
$apiDef = @(": a collection of exposed interfaces and protocols for the purpose of general reuse")
$null = Insert-TermBefore $d $constantPara "Application Program Interface (API)" $apiDef

# ---------------------------------------------------------------------------
# 2) "Comment" and "Documentation" and "Integrated Development Environment
#    (IDE)" after "Constant"
# ---------------------------------------------------------------------------
$constantPara = Find-Paragraph $d "Constant:"
$commentDef = @(": an annotation of a line or section of code")
$commentPara = Insert-TermAfter $d $constantPara "Comment" $commentDef

$documentationDef = @(": ", "literature that provides details about a library or tool")
$documentationPara = Insert-TermAfter $d $commentPara "Documentation" $documentationDef

$idePara = Insert-TermAfter $d $documentationPara "Integrated Development Environment (IDE)" @(": a source code editor with automation tools")
Add-BookmarkAtParagraphEnd $d $idePara "_GoBack"

# ---------------------------------------------------------------------------
# 3) "Library" before "Loop"
# ---------------------------------------------------------------------------
$loopPara = Find-Paragraph $d "Loop:"
$libraryDef = @(": a collection of implementations for the purpose of general reuse")
$null = Insert-TermBefore $d $loopPara "Library" $libraryDef

# ---------------------------------------------------------------------------
# 4) "Pseudo Code" and "Readability" before "Recursion"
# ---------------------------------------------------------------------------
$recursionPara = Find-Paragraph $d "Recursion:"
$pseudoDef = @(": an informal or simplified programming language used to describe how a program should execute")
$pseudoPara = Insert-TermBefore $d $recursionPara "Pseudo Code" $pseudoDef

$readabilityDef = @(": the measure of ease of interpretation and understanding of source code")
$null = Insert-TermAfter $d $pseudoPara "Readability" $readabilityDef

Write-Host "Done"
